# ---------------------------------------------------------------------------
# Enabled parallel cross-browser testing to optimize execution time and coverage
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# --- Sheet "Tests": mark the "loginTestWithInValidCredentials" row to execute ---
$wsTests = $wb.Worksheets.Item("Tests")
$wsTests.Range("C3").Value = "yes"

# --- Sheet "DataProviderTests": add a "browser" column and cross-browser rows ---
$wsData = $wb.Worksheets.Item("DataProviderTests")

# Insert a new column before the current "username" column (col C) so that the
# layout becomes: testname | execute | browser | username | password
$wsData.Columns("C").Insert()

# Header row
$wsData.Range("C1").Value = "browser"

# Match the formatting already used by the other header cells (bold, centered)
$wsData.Range("C1").HorizontalAlignment = -4108
$wsData.Range("C1").VerticalAlignment = -4108
$wsData.Range("C1").Font.Bold = $true

# Rebuild the full data block (rows 2-9) for testname / execute / browser / username / password
$data = @(
    @("loginTestWithValidCredentials",   "yes", "chrome",  "Admin",    "admin123"),
    @("loginTestWithValidCredentials",   "yes", "edge",    "Admin",    "admin123"),
    @("loginTestWithValidCredentials",   "yes", "firefox", "Admin",    "admin123"),
    @("loginTestWithInValidCredentials", "yes", "chrome",  "Admin",    "admin12345"),
    @("loginTestWithValidCredentials",   "no",  "",        "Admin",    "admin123"),
    @("loginTestWithValidCredentials",   "no",  "",        "adminnn",  "admin123"),
    @("loginTestWithInValidCredentials", "yes", "edge",    "Admin",    "admin12345"),
    @("loginTestWithInValidCredentials", "yes", "firefox", "Admin",    "admin12345")
)

# Populate the non-browser columns first (A, B, D, E) for every row ...
$row = 2
foreach ($r in $data) {
    $wsData.Range("A$row").Value = $r[0]
    $wsData.Range("B$row").Value = $r[1]
    $wsData.Range("D$row").Value = $r[3]
    $wsData.Range("E$row").Value = $r[4]
    $row = $row + 1
}

# ... then write the "browser" column in "chrome, firefox, edge" first-seen order
# (row 2 = chrome, row 4 = firefox, row 3 = edge, remaining rows reuse those values)
$wsData.Range("C2").Value = "chrome"
$wsData.Range("C4").Value = "firefox"
$wsData.Range("C3").Value = "edge"
$wsData.Range("C5").Value = "chrome"
$wsData.Range("C8").Value = "edge"
$wsData.Range("C9").Value = "firefox"

foreach ($rn in 2,3,4,5,8,9) {
    $wsData.Range("C$rn").HorizontalAlignment = -4108
    $wsData.Range("C$rn").VerticalAlignment = -4108
}

# Rows 6 and 7 ("no"-execute rows) have no browser assigned - leave that cell blank
$wsData.Range("C6").Clear()
$wsData.Range("C7").Clear()

# The existing "username"/"password" columns keep their width automatically as they
# shift right with the insert; only the new "browser" column needs an explicit width.
$wsData.Columns("C").ColumnWidth = 7.5

$wsTests.Range("D9").Select()
$wsData.Activate()
$wsData.Range("F7").Select()
